$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Lhcgr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.052798
$ws.Range("N2").Value = 0.158394
$ws.Range("O2").Value = 0.1385348856125547
$ws.Range("P2").Value = 0.1385348856125547
$ws.Range("Q2").Value = 8.971206858286001
$ws.Range("R2").Value = 80.740861724574
$ws.Range("S2").Value = 0.06152590465020416
$ws.Range("T2").Value = 0.06152590465020417

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Lhcgr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.328319
$ws.Range("N3").Value = 0.9849570000000001
$ws.Range("O3").Value = 0.8614651143874453
$ws.Range("P3").Value = 0.8614651143874454
$ws.Range("Q3").Value = 55.78653859058301
$ws.Range("R3").Value = 502.0788473152471
$ws.Range("S3").Value = 0.3825925885232467
$ws.Range("T3").Value = 0.3825925885232468

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Lhcgr"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 68.382243
$ws.Range("H4").Value = 205.146729
$ws.Range("I4").Value = 0.1787346690539575
$ws.Range("J4").Value = 0.1787346690539575
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.052798
$ws.Range("N4").Value = 0.158394
$ws.Range("O4").Value = 0.1385348856125547
$ws.Range("P4").Value = 0.1385348856125547
$ws.Range("Q4").Value = 3.610445665914
$ws.Range("R4").Value = 32.494010993226
$ws.Range("S4").Value = 0.02476098693238782
$ws.Range("T4").Value = 0.02476098693238782

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Lhcgr"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 68.382243
$ws.Range("H5").Value = 205.146729
$ws.Range("I5").Value = 0.1787346690539575
$ws.Range("J5").Value = 0.1787346690539575
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.328319
$ws.Range("N5").Value = 0.9849570000000001
$ws.Range("O5").Value = 0.8614651143874453
$ws.Range("P5").Value = 0.8614651143874454
$ws.Range("Q5").Value = 22.451189639517
$ws.Range("R5").Value = 202.060706755653
$ws.Range("S5").Value = 0.1539736821215697
$ws.Range("T5").Value = 0.1539736821215697

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Lhcgr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 53.27463399999999
$ws.Range("H6").Value = 159.823902
$ws.Range("I6").Value = 0.1392470275793777
$ws.Range("J6").Value = 0.1392470275793778
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.052798
$ws.Range("N6").Value = 0.158394
$ws.Range("O6").Value = 0.1385348856125547
$ws.Range("P6").Value = 0.1385348856125547
$ws.Range("Q6").Value = 2.812794125932
$ws.Range("R6").Value = 25.315147133388
$ws.Range("S6").Value = 0.01929057103759734
$ws.Range("T6").Value = 0.01929057103759734

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Lhcgr"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 53.27463399999999
$ws.Range("H7").Value = 159.823902
$ws.Range("I7").Value = 0.1392470275793777
$ws.Range("J7").Value = 0.1392470275793778
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.328319
$ws.Range("N7").Value = 0.9849570000000001
$ws.Range("O7").Value = 0.8614651143874453
$ws.Range("P7").Value = 0.8614651143874454
$ws.Range("Q7").Value = 17.491074560246
$ws.Range("R7").Value = 157.419671042214
$ws.Range("S7").Value = 0.1199564565417804
$ws.Range("T7").Value = 0.1199564565417804

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Lhcgr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 91.01828266666666
$ws.Range("H8").Value = 273.054848
$ws.Range("I8").Value = 0.2378998101932138
$ws.Range("J8").Value = 0.2378998101932138
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.052798
$ws.Range("N8").Value = 0.158394
$ws.Range("O8").Value = 0.1385348856125547
$ws.Range("P8").Value = 0.1385348856125547
$ws.Range("Q8").Value = 4.805583288234667
$ws.Range("R8").Value = 43.250249594112
$ws.Range("S8").Value = 0.03295742299236534
$ws.Range("T8").Value = 0.03295742299236534

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Lhcgr"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 91.01828266666666
$ws.Range("H9").Value = 273.054848
$ws.Range("I9").Value = 0.2378998101932138
$ws.Range("J9").Value = 0.2378998101932138
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.328319
$ws.Range("N9").Value = 0.9849570000000001
$ws.Range("O9").Value = 0.8614651143874453
$ws.Range("P9").Value = 0.8614651143874454
$ws.Range("Q9").Value = 29.88303154683733
$ws.Range("R9").Value = 268.947283921536
$ws.Range("S9").Value = 0.2049423872008485
$ws.Range("T9").Value = 0.2049423872008485
